$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '30.489.38'
Set-TextValue "E2" '  +1.03%  '
Set-TextValue "D3" '1.879.78'
Set-TextValue "E3" '  +0.92%  '
Set-TextValue "D4" '0.9998'
Set-TextValue "E4" '  -0.13%  '
Set-TextValue "D5" '246.52'
Set-TextValue "E5" '  +5.55%  '
Set-TextValue "D6" '0.9998'
Set-TextValue "E6" '  -0.11%  '
Set-TextValue "D7" '0.4763'
Set-TextValue "E7" '  +2.08%  '
Set-TextValue "D8" '0.2903'
Set-TextValue "E8" '  +2.05%  '
Set-TextValue "D9" '0.06519'
Set-TextValue "E9" '  +0.71%  '
Set-TextValue "D10" '21.95'
Set-TextValue "E10" '  +3.27%  '
Set-TextValue "D11" '0.07745'
Set-TextValue "E11" '  +0.13%  '
Set-TextValue "D12" '0.7433'
Set-TextValue "E12" '  +8.95%  '
Set-TextValue "D13" '97.02'
Set-TextValue "E13" '  +3.58%  '
Set-TextValue "D14" '1.872.71'
Set-TextValue "E14" '  +0.57%  '
Set-TextValue "D15" '5.131'
Set-TextValue "E15" '  +1.95%  '
Set-TextValue "D16" '274.36'
Set-TextValue "E16" '  +1.91%  '
Set-TextValue "D17" '30.475.18'
Set-TextValue "E17" '  +1.02%  '
Set-TextValue "D18" '13.64'
Set-TextValue "E18" '  +2.50%  '
Set-TextValue "D19" '0.000007581'
Set-TextValue "E19" '  -0.47%  '
Set-TextValue "D21" '2.131.69'
Set-TextValue "E21" '  +1.78%  '
Set-TextValue "D22" '0.9999'
Set-TextValue "E22" '  -0.13%  '
Set-TextValue "D24" '6.176'
Set-TextValue "E24" '  +1.40%  '
Set-TextValue "D25" '9.290'
Set-TextValue "E25" '  -0.41%  '
Set-TextValue "D26" '164.48'
Set-TextValue "E26" '  -0.58%  '
Set-TextValue "D27" '18.93'
Set-TextValue "E27" '  +2.12%  '
Set-TextValue "D28" '1.961'
Set-TextValue "E28" '  +3.97%  '
Set-TextValue "D30" '0.09994'
Set-TextValue "E30" '  +1.93%  '
Set-TextValue "D31" '1.514'
Set-TextValue "E31" '  +4.37%  '
Set-TextValue "D32" '4.332'
Set-TextValue "E32" '  +2.47%  '
Set-TextValue "D34" '0.04785'
Set-TextValue "E34" '  +2.56%  '
Set-TextValue "D35" '1.126'
Set-TextValue "E35" '  +0.70%  '
Set-TextValue "D36" '0.6986'
Set-TextValue "E36" '  +1.39%  '
Set-TextValue "D37" '2.715'
Set-TextValue "E37" '  +0.17%  '
Set-TextValue "D38" '0.01868'
Set-TextValue "E38" '  +2.00%  '
Set-TextValue "D39" '2.730'
Set-TextValue "E39" '  -0.60%  '
Set-TextValue "D40" '6.354'
Set-TextValue "E40" '  +0.98%  '
Set-TextValue "D41" '1.938'
Set-TextValue "E41" '  +3.06%  '
Set-TextValue "D42" '70.05'
Set-TextValue "E42" '  -0.97%  '
Set-TextValue "D43" '0.4178'
Set-TextValue "E43" '  +3.18%  '
Set-TextValue "D44" '0.9996'
Set-TextValue "E44" '  -0.07%  '
Set-TextValue "D46" '102.87'
Set-TextValue "E46" '  +0.81%  '
Set-TextValue "D47" '9.337'
Set-TextValue "E47" '  +3.46%  '
Set-TextValue "D48" '7.099'
Set-TextValue "E48" '  +2.31%  '
Set-TextValue "D49" '35.37'
Set-TextValue "E49" '  +3.99%  '
Set-TextValue "D50" '922.32'
Set-TextValue "E50" '  -1.30%  '
Set-TextValue "D51" '0.05612'
Set-TextValue "E51" '  +0.69%  '
Set-TextValue "D33" '4.061'
Set-TextValue "E23" '  +2.25%  '
Set-TextValue "E29" '  +0.55%  '
Set-TextValue "E45" '  +0.75%  '
